$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style used by
# the existing header cells (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for the new I and J columns, one entry per data row (rows 2-32).
$iVals = @(8, 7, 5, 6, 9, 8, 9, 9, 8, 8, 7, 8, 8, 9, 4, 8, 7, 6, 7, 8, 6, 8, 6, 8, 7, 7, 6, 8, 6, 5, 5)
$jVals = @(8, 8, 5, 6, 9, 8, 9, 9, 8, 8, 7, 8, 8, 9, 4, 8, 7, 6, 7, 8, 6, 8, 6, 8, 7, 7, 6, 8, 6, 5, 5)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}

Write-Output "I0 and IF columns added"
